$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert a comma after "As an administrator" in every admin user-story
#    paragraph: "As an administrator I would" -> "As an administrator, I would"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "As an administrator I would", $false, $false, $false, $false, $false,
    $true, 1, $false, "As an administrator, I would", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. In the "safe locations" user story, capitalise the stray "i" that starts
#    the second sentence: "... general location then i would also ..." ->
#    "... general location then I would also ..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "general location then i would also", $false, $false, $false, $false,
    $false, $true, 1, $false, "general location then I would also", 2) | Out-Null

# ---------------------------------------------------------------------------
# Helper: find a target paragraph by its exact text and return its 1-based
# index within $d.Paragraphs. (Going through Paragraphs.Item() on a Range
# returned from Find is unreliable, so we locate the index by counting the
# paragraphs that precede the match instead.)
# ---------------------------------------------------------------------------
function Get-ParagraphIndexForText($searchText) {
    $probe = $d.Content
    $found = $probe.Find.Execute($searchText, $false, $false, $false, $false,
                                  $false, $true, 1, $false, "", 0)
    if (-not $found) {
        return -1
    }
    $pre = $d.Range(0, $probe.Start)
    return $pre.Paragraphs.Count + 1
}

# ---------------------------------------------------------------------------
# 3. Add the new "contact details" user story right before the existing
#    "quickly see some information about my account" story - once for the
#    administrator section and once for the user section.
# ---------------------------------------------------------------------------
$adminIdx = Get-ParagraphIndexForText("As an administrator, I would like to be able to quickly see some information about my account.")
if ($adminIdx -gt 0) {
    $target = $d.Paragraphs.Item($adminIdx)
    $target.Range.InsertParagraphBefore()
    $newPara = $d.Paragraphs.Item($adminIdx)
    $newPara.Range.Text = "As an administrator, I would like to see contact details about the platform."
}

$userIdx = Get-ParagraphIndexForText("As a user, I would like to be able to quickly see some information about my account.")
if ($userIdx -gt 0) {
    $target2 = $d.Paragraphs.Item($userIdx)
    $target2.Range.InsertParagraphBefore()
    $newPara2 = $d.Paragraphs.Item($userIdx)
    $newPara2.Range.Text = "As a user, I would like to see contact details about the platform."
}

# ---------------------------------------------------------------------------
# 4. Re-write the "User Stories Users" heading so the stale
#    lastRenderedPageBreak rendering cache is dropped (the run is otherwise
#    unchanged).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "User Stories Users", $false, $false, $false, $false, $false, $true, 1,
    $false, "User Stories Users", 2) | Out-Null

Write-Output "done"
